$d = $word.ActiveDocument

# 1) Update the letter date: "September 19, 2025" -> "September 21, 2025".
#    The paragraph wraps its run in a Text10 bookmark, so rewrite the whole
#    paragraph via InsertXML (keeping the bookmark + run properties intact)
#    rather than a plain Find/Replace, which would silently drop the
#    xml:space="preserve" attribute on the text run.
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -eq "September 19, 2025`r") {
        $dateXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr><w:bookmarkStart w:id="1" w:name="Text10"/><w:bookmarkEnd w:id="1"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve">September 21, 2025</w:t></w:r></w:p>'
        [void]$para.Range.InsertXML($dateXml)
        break
    }
}

# 2) Split the mailing-address paragraph "969 Story Road, San Jose CA 95122"
#    into two paragraphs: "969 Story Road" and "San Jose, CA 95122".
#    Only the first occurrence (the letter's return-address block) changes;
#    the identical text inside the property-address table must stay intact,
#    so we scope the edit to that specific paragraph via the Paragraphs
#    collection, then rewrite it with InsertXML for exact formatting control.
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -eq "969 Story Road, San Jose CA 95122`r") {
        $addrXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve">969 Story Road</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve">San Jose, CA 95122</w:t></w:r></w:p>'
        [void]$para.Range.InsertXML($addrXml)
        break
    }
}

# 3) Remove the blank "NoSpacing" paragraph that immediately follows the
#    "...Board of Directors" signature line.
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -match "Board of Directors") {
        $blank = $d.Paragraphs.Item($i + 1)
        if ($blank.Range.Text -eq "`r") {
            [void]$blank.Range.Delete()
        }
        break
    }
}
